# Update crypto price (D) and 1h volume-change (E) columns to the latest
# scraped values, as produced by the GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.457.48'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.607.72'
$ws.Range("E3").Value = '  -0.68%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.53'
$ws.Range("E5").Value = '  -1.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '190.25'
$ws.Range("E6").Value = '  -0.97%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.604.78'
$ws.Range("E7").Value = '  -0.65%  '
$ws.Range("E8").Value = '  -1.97%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.185'
$ws.Range("E10").Value = '  +4.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.663'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '56.28'
$ws.Range("E12").Value = '  -3.60%  '
$ws.Range("E13").Value = '  +7.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.73'
$ws.Range("E14").Value = '  -2.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.186.48'
$ws.Range("E15").Value = '  -0.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.85'
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.604.21'
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.367.47'
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.67'
$ws.Range("E19").Value = '  -0.20%  '
$ws.Range("E20").Value = '  +0.32%  '
$ws.Range("E21").Value = '  -0.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '491.88'
$ws.Range("E22").Value = '  +0.68%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.38'
$ws.Range("E23").Value = '  -0.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.90'
$ws.Range("E24").Value = '  -8.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.22'
$ws.Range("E25").Value = '  +6.71%  '
$ws.Range("E26").Value = '  -1.94%  '
$ws.Range("E27").Value = '  -4.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.07'
$ws.Range("E28").Value = '  -2.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.44'
$ws.Range("E29").Value = '  -2.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.31'
$ws.Range("E30").Value = '  -2.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.58'
$ws.Range("E31").Value = '  -3.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.29'
$ws.Range("E32").Value = '  -0.34%  '
$ws.Range("E33").Value = '  -1.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '66.01'
$ws.Range("E34").Value = '  -0.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '581.69'
$ws.Range("E35").Value = '  -7.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '39.06'
$ws.Range("E36").Value = '  +0.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0₃0819'
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("E39").Value = '  -3.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.98'
$ws.Range("E40").Value = '  +6.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.28'
$ws.Range("E41").Value = '  +20.24%  '
$ws.Range("E42").Value = '  -2.89%  '
$ws.Range("E43").Value = '  -7.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.227.35'
$ws.Range("E44").Value = '  -2.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.06'
$ws.Range("E45").Value = '  -1.68%  '
$ws.Range("E46").Value = '  -1.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.71'
$ws.Range("E47").Value = '  +6.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.39'
$ws.Range("E48").Value = '  +2.61%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("E50").Value = '  -2.16%  '
$ws.Range("E51").Value = '  +0.09%  '
